# Update the "Förändrad" (Changed) date column C for rows 2-39 from
# 2023-09-03 (45172) to 2023-09-06 (45175).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C39").Value = 45175
